$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename worksheet tab: "Deal Success" -> "Deal" ---
$ws.Name = "Deal"

# --- Header row (row 1) text updates ---
# C1 header: navpeII__Company_Name__c -> navpeII__Company__c
$ws.Range("C1").Value = "navpeII__Company__c"

# --- New "Status" column A, filled with "Created" for every data row ---
$ws.Range("A2:A6").Value = "Created"

# --- Stage column (D) value updates ---
$ws.Range("D3").Value = "Idea Generation"
$ws.Range("D4").Value = "New/Initial Interest"
$ws.Range("D5").Value = "New/Initial Interest"
$ws.Range("D6").Value = "New/Initial Interest"

# --- Log-in-date column (E) : rows 5 & 6 were real dates, convert them to the
#     same literal text "08/31/2022" used elsewhere in the column ---
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "08/31/2022"
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "08/31/2022"
$ws.Range("E6").Style = "Normal"

# --- Header font formatting ---
# B1, D1, E1 -> Segoe UI 10, color #181818
foreach ($addr in @("B1", "D1", "E1")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Segoe UI"
    $c.Font.Size = 10
    $c.Font.Color = 1579032
}

# C1 -> Segoe UI 10, color #181818, wrap text, left/center aligned
$ws.Range("C1").Font.Name = "Segoe UI"
$ws.Range("C1").Font.Size = 10
$ws.Range("C1").Font.Color = 1579032
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("C1").WrapText = $true

# --- Stage column data font formatting (D2:D6) -> Arial 9, black ---
$stage = $ws.Range("D2:D6")
$stage.Font.Name = "Arial"
$stage.Font.Size = 9
$stage.Font.Color = 0

# --- Row / column sizing ---
$ws.Rows.Item(1).RowHeight = 42.75
$ws.Columns.Item(4).ColumnWidth = 14.6
$ws.Columns.Item(5).ColumnWidth = 22.76

# --- Selection state ---
$ws.Range("E2").Select() | Out-Null

# --- Page setup ---
$ws.PageSetup.Orientation = 1 | Out-Null
